$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 118, shifting existing rows 118:223 down to 119:224
$ws.Rows.Item(118).Insert()

# Populate the newly inserted row 118 with a new weekly price record
# (same Mercado/Region/Categoria/Unidad/Origen/Clasificacion as the template row, new date & price data)
$ws.Cells.Item(118, 1).Value = 4
$ws.Cells.Item(118, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(118, 3).Value = "Los Lagos"
$ws.Cells.Item(118, 4).Value = 44586
$ws.Cells.Item(118, 5).Value = 10
$ws.Cells.Item(118, 6).Value = 100112040
$ws.Cells.Item(118, 7).Value = "Cilantro"
$ws.Cells.Item(118, 8).Value = "Sin especificar"
$ws.Cells.Item(118, 9).Value = "Primera"
$ws.Cells.Item(118, 10).Value = 80
$ws.Cells.Item(118, 11).Value = 10000
$ws.Cells.Item(118, 12).Value = 10000
$ws.Cells.Item(118, 13).Value = 10000
$ws.Cells.Item(118, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(118, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(118, 16).Value = 5000
$ws.Cells.Item(118, 17).Value = 2
$ws.Cells.Item(118, 18).Value = "Hortaliza"

# Match the date number format used by the rest of column D
$ws.Cells.Item(118, 4).NumberFormat = $ws.Cells.Item(119, 4).NumberFormat
